# Updates cryptos list price/volume figures (and two pairs of swapped
# coin rows) to match the latest scrape, per the commit diff.
#
# Price-column values are written with a leading apostrophe where the
# text looks like a plain number (e.g. "0.9982") so Excel keeps storing
# them as text (matching the original inlineStr cells) instead of
# silently converting them to numeric cells. Values that already contain
# multiple dots (e.g. "28.234.47") are never auto-converted by Excel, so
# they are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "28.234.47"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.810.70"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = "  -0.34%  "

# Row 5: BNB
$ws.Range("D5").Value = "'317.08"
$ws.Range("E5").Value = "  -0.28%  "

# Row 6: USDC
$ws.Range("D6").Value = "'0.9987"
$ws.Range("E6").Value = "  -0.28%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.5325"
$ws.Range("E7").Value = "  -1.94%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.3970"
$ws.Range("E8").Value = "  +4.68%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.07551"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10: OKB
$ws.Range("D10").Value = "'41.64"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11: Polygon
$ws.Range("D11").Value = "'1.095"
$ws.Range("E11").Value = "  -1.78%  "

# Row 12: BinanceUSD
$ws.Range("D12").Value = "'0.9977"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13: Chainlink
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'7.591"
$ws.Range("E13").Value = "  +3.23%  "

# Row 14: Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.258"
$ws.Range("E14").Value = "  +1.57%  "

# Row 15: Solana
$ws.Range("D15").Value = "'20.54"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "1.808.48"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17: Litecoin
$ws.Range("D17").Value = "'89.16"
$ws.Range("E17").Value = "  -1.18%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "'0.00001066"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19: TRON
$ws.Range("D19").Value = "'0.06585"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20: Dai
$ws.Range("D20").Value = "'0.9988"
$ws.Range("E20").Value = "  -0.20%  "

# Row 21: Avalanche
$ws.Range("D21").Value = "'17.35"
$ws.Range("E21").Value = "  -0.50%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'6.022"
$ws.Range("E22").Value = "  +1.16%  "

# Row 23: WrappedBTC
$ws.Range("D23").Value = "28.226.39"
$ws.Range("E23").Value = "  -1.53%  "

# Row 24: Cosmos
$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25: Toncoin
$ws.Range("D25").Value = "'2.089"
$ws.Range("E25").Value = "  +0.57%  "

# Row 26: Monero
$ws.Range("D26").Value = "'156.44"
$ws.Range("E26").Value = "  -2.96%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'20.35"
$ws.Range("E27").Value = "  -0.81%  "

# Row 28: WrappedliquidstakedEther2.0
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.012.04"
$ws.Range("E28").Value = "  +0.14%  "

# Row 29: LidoDAOToken
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.375"
$ws.Range("E29").Value = "  +1.44%  "

# Row 30: BitcoinCash
$ws.Range("D30").Value = "'123.10"
$ws.Range("E30").Value = "  -0.09%  "

# Row 31: Stellar
$ws.Range("D31").Value = "'0.1098"
$ws.Range("E31").Value = "  +3.62%  "

# Row 32: ImmutableX
$ws.Range("D32").Value = "'1.102"
$ws.Range("E32").Value = "  -2.71%  "

# Row 33: HuobiToken
$ws.Range("D33").Value = "'3.670"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34: Filecoin
$ws.Range("D34").Value = "'5.555"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35: Hedera
$ws.Range("D35").Value = "'0.07178"
$ws.Range("E35").Value = "  +8.52%  "

# Row 36: Algorand
$ws.Range("D36").Value = "'0.2230"
$ws.Range("E36").Value = "  -1.45%  "

# Row 37: InternetComputer(DFINITY)
$ws.Range("D37").Value = "'5.209"
$ws.Range("E37").Value = "  +3.00%  "

# Row 38: VeChain
$ws.Range("D38").Value = "'0.02304"
$ws.Range("E38").Value = "  -0.45%  "

# Row 39: FraxShare
$ws.Range("D39").Value = "'8.615"
$ws.Range("E39").Value = "  -0.05%  "

# Row 40: Aptos
$ws.Range("D40").Value = "'11.28"
$ws.Range("E40").Value = "  +0.18%  "

# Row 41: TheSandbox
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6191"
$ws.Range("E41").Value = "  -0.76%  "

# Row 42: TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.195"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43: WEMIXTOKEN
$ws.Range("D43").Value = "'1.404"
$ws.Range("E43").Value = "  -3.22%  "

# Row 44: EnergySwap
$ws.Range("D44").Value = "'13.46"

# Row 45: PancakeSwap
$ws.Range("D45").Value = "'3.687"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46: Decentraland
$ws.Range("D46").Value = "'0.5767"
$ws.Range("E46").Value = "  -1.53%  "

# Row 47: Quant
$ws.Range("D47").Value = "'125.11"
$ws.Range("E47").Value = "  -1.74%  "

# Row 48: NEARProtocol
$ws.Range("D48").Value = "'1.944"
$ws.Range("E48").Value = "  -1.04%  "

# Row 49: EOS
$ws.Range("E49").Value = "  +0.28%  "

# Row 50: Cronos
$ws.Range("D50").Value = "'0.06826"
$ws.Range("E50").Value = "  -1.20%  "

# Row 51: Aave
$ws.Range("D51").Value = "'71.23"
$ws.Range("E51").Value = "  -1.89%  "
